$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.349.13"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.95%  "

$ws.Range("D3").Value = "'3.751.34"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.61%  "

$ws.Range("D5").Value = "'602.58"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("D6").Value = "'168.32"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.51%  "

$ws.Range("D7").Value = "'3.749.99"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.63%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.17%  "

$ws.Range("E9").Value = "  +1.20%  "

$ws.Range("E10").Value = "  +4.39%  "

$ws.Range("E11").Value = "  -0.17%  "

$ws.Range("E12").Value = "  +0.89%  "

$ws.Range("D13").Value = "'38.28"
$ws.Range("D13").ClearFormats()

$ws.Range("E14").Value = "  +1.96%  "

$ws.Range("D15").Value = "'4.378.33"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.54%  "

$ws.Range("D16").Value = "'3.747.21"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.44%  "

$ws.Range("D17").Value = "'69.302.33"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.89%  "

$ws.Range("D18").Value = "'7.43"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.64%  "

$ws.Range("E19").Value = "  +1.02%  "

$ws.Range("E20").Value = "  -1.42%  "

$ws.Range("E21").Value = "  +12.38%  "

$ws.Range("D22").Value = "'493.15"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.39%  "

$ws.Range("E23").Value = "  +1.19%  "

$ws.Range("E24").Value = "  +5.41%  "

$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("E26").Value = "  -0.14%  "

$ws.Range("D27").Value = "'12.33"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.20%  "

$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("E30").Value = "  +1.71%  "

$ws.Range("D31").Value = "'8.17"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.47%  "

$ws.Range("D32").Value = "'2.46"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.09%  "

$ws.Range("D33").Value = "'31.69"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.28%  "

$ws.Range("D34").Value = "'3.895.71"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.49%  "

$ws.Range("E35").Value = "  +0.98%  "

$ws.Range("D36").Value = "'3.683.98"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.56%  "

$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.140"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.90%  "

$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").Value = "'5.99"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.31%  "

$ws.Range("E40").Value = "  +0.34%  "

$ws.Range("E41").Value = "  +1.08%  "

$ws.Range("D42").Value = "'3.07"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.97%  "

$ws.Range("D43").Value = "'48.94"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.44%  "

$ws.Range("E44").Value = "  +1.51%  "

$ws.Range("D45").Value = "'424.36"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.18%  "

$ws.Range("E46").Value = "  +0.85%  "

$ws.Range("D48").Value = "'40.26"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.99%  "

$ws.Range("D49").Value = "'141.21"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.23%  "

$ws.Range("D50").Value = "'2.797.02"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.08%  "

$ws.Range("D51").Value = "'0.0356"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.88%  "

